$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.709.59"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.637.51"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Formula = "=""606.97"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Formula = "=""146.98"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Formula = "=""0.588"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("D10").Formula = "=""0.382"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +5.97%  "
$ws.Range("D11").Formula = "=""5.60"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "3.110.86"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "63.564.67"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "2.639.19"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Formula = "=""11.75"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Formula = "=""4.55"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("D20").Formula = "=""346.13"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Formula = "=""6.89"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Formula = "=""5.57"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").Formula = "=""66.25"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("E25").Value = "  +7.31%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").Formula = "=""9.21"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +6.23%  "
$ws.Range("D28").Formula = "=""567.07"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("D29").Formula = "=""8.07"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "0.0₃0851"
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("D36").Formula = "=""168.61"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Formula = "=""1.95"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +5.46%  "
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Formula = "=""164.74"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -6.89%  "
$ws.Range("D43").Formula = "=""40.14"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Formula = "=""3.78"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").Formula = "=""21.84"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("D48").Formula = "=""2.00"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +14.10%  "
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Formula = "=""0.0955"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").Formula = "=""18.74"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -1.63%  "
